$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2020" column (column K), mirroring the
# formatting of the existing last data column (J) for each affected row.

# Row 3: thin bottom-border-only separator cell, no value.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)  # xlPasteFormats

# Row 4: header year value.
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(4, 11).Value = 2020

# Row 5: data value.
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(5, 11).Value = 173

# Match the saved selection/active cell recorded in the sheet view.
$ws.Range("I18").Select()
